# Update the "current attendees" counts (column F) on the "展览" and
# "全部类型" sheets. Both sheets carry the same rows, so the same set of
# updates is applied twice.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1087
    4  = 362
    5  = 4650
    7  = 392
    8  = 1383
    9  = 919
    11 = 1119
    13 = 612
    15 = 25
    16 = 9
    17 = 273
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
